$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mon Sep 11 14:06:11 EDT 2023"
$ws.Range("B3").Value = "Mon Sep 11 14:06:24 EDT 2023"
$ws.Range("B4").Value = "Mon Sep 11 14:06:38 EDT 2023"
